$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Luke": update rows 5-6 dates, fill in rows 7-9 with new task data
# ----------------------------------------------------------------------
$luke = $wb.Worksheets.Item("Luke")

# Rows 5 & 6: date moved from 42800 (6 Mar 2017) to 42861 (6 May 2017)
$luke.Range("A5").Value = 42861
$luke.Range("A6").Value = 42861

# New shared-string text is entered in this precise order so new entries
# land at shared-string indices 32-39 in the same order the original
# author typed them in:
#   32 Login Section, 33 Starting the Modal login, 34 T05, 35 T04,
#   36 Links, 37 Establishing correct links on home page,
#   38 9/05.2017, 39 Completed Modal Login section without DB
$luke.Range("E7").Value = "Login Section"
$luke.Range("F7").Value = "Starting the Modal login"
$luke.Range("D7").Value = "T05"
$luke.Range("D8").Value = "T04"
$luke.Range("E8").Value = "Links"
$luke.Range("F8").Value = "Establishing correct links on home page"
$luke.Range("A7").Value = "9/05.2017"
$luke.Range("F9").Value = "Completed Modal Login section without DB"

# Remaining cells in rows 7-9 (numbers + re-used shared strings)
$luke.Range("B7").Value = 3
$luke.Range("C7").Value = "S27"

# Row 8: real date 42865 (10 May 2017) - copy date style from A6 first so the
# cell reuses the existing date cellXfs entry instead of creating a new one
$luke.Range("A6").Copy($luke.Range("A8"))
$luke.Range("A8").Value = 42865
$luke.Range("B8").Value = 2
$luke.Range("C8").Value = "S27"

# Row 9: real date 42865 (10 May 2017)
$luke.Range("A6").Copy($luke.Range("A9"))
$luke.Range("A9").Value = 42865
$luke.Range("B9").Value = 4
$luke.Range("C9").Value = "S27"
$luke.Range("D9").Value = "T05"
$luke.Range("E9").Value = "Login Section"

# ----------------------------------------------------------------------
# Sheet "BurnDown": log hours on rows 28 & 29 (updates the H "Hours
# Remaining" running formula and the G59 total automatically)
# ----------------------------------------------------------------------
$burn = $wb.Worksheets.Item("BurnDown")
$burn.Range("G28").Value = 3
$burn.Range("G29").Value = 6

# ----------------------------------------------------------------------
# Selections: BurnDown keeps its own (non-active) selection at G29; Luke
# (the active tab) ends up selected at B11
# ----------------------------------------------------------------------
$burn.Range("G29").Select()
$luke.Range("B11").Select()

Write-Host "done"
